$d = $word.ActiveDocument

$replacements = @(
    @{ Old = "13×32="; New = "89×18=" },
    @{ Old = "75×72="; New = "33×89=" },
    @{ Old = "69×38="; New = "41×65=" },
    @{ Old = "97×76="; New = "50×80=" },
    @{ Old = "74×79="; New = "75×49=" },
    @{ Old = "72×36="; New = "99×80=" },
    @{ Old = "82×88="; New = "57×16=" },
    @{ Old = "28×54="; New = "88×20=" },
    @{ Old = "82×49="; New = "41×25=" },
    @{ Old = "82×32="; New = "16×61=" },
    @{ Old = "12×19="; New = "80×87=" },
    @{ Old = "17×48="; New = "29×45=" },
    @{ Old = "82×35="; New = "83×48=" },
    @{ Old = "77×15="; New = "20×85=" },
    @{ Old = "56×19="; New = "76×92=" },
    @{ Old = "56×59="; New = "22×45=" },
    @{ Old = "63×38="; New = "62×69=" },
    @{ Old = "17×98="; New = "18×25=" },
    @{ Old = "33×15="; New = "73×46=" },
    @{ Old = "21×38="; New = "35×65=" },
    @{ Old = "25×51="; New = "18×52=" },
    @{ Old = "78×60="; New = "70×24=" },
    @{ Old = "31×30="; New = "32×94=" },
    @{ Old = "56×84="; New = "42×69=" },
    @{ Old = "25×80="; New = "63×58=" }
)

foreach ($r in $replacements) {
    $range = $d.Content
    $range.Find.Execute($r.Old, $true, $true, $false, $false, $false,
                         $true, 1, $false, $r.New, 2)
}
